$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows after the header row, pushing the existing data rows
# (currently rows 2-19) down to rows 5-22.
$ws.Rows("2:4").Insert()

# Copy the formatting (cell borders/style) of the first data row (now row 5,
# which used to be row 2) onto the freshly inserted blank rows so they match
# the look of the rest of the table.
$ws.Range("A5:C5").Copy()
$ws.Range("A2:C4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new "bathrooms" / "bathroom_text" / "beds" comment rows.
# B2 is intentionally left blank (no "Replace null with" value for bathroom_text).
$ws.Range("A2").Value = "bathroom_text"
$ws.Range("C2").Value = "Turned text string to float and saved float values in ""bathroom"" colummn"

$ws.Range("A3").Value = "bathrooms"
$ws.Range("B3").Value = "mean"
$ws.Range("C3").Value = "Used cleaned values from bathroom text column, used mean for null values, outliers seemed reasonable "

$ws.Range("A4").Value = "beds"
$ws.Range("B4").Value = "mean"
$ws.Range("C4").Value = "Kept it simple and used mean, ouliers seemed reasonable"

# Widen column C slightly, matching the author's manual resize.
$ws.Columns("C").ColumnWidth = 87.66666666666667

# Match the saved selection/active cell.
$ws.Range("C6").Select()
